$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 119; this pushes existing rows 119:248 down to 120:249
$ws.Rows.Item(119).Insert()

# Populate the newly inserted row 119 with the new data record
$ws.Cells.Item(119, 1).Value = 3
$ws.Cells.Item(119, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(119, 3).Value = "Coquimbo"
$ws.Cells.Item(119, 4).Value = 44944
$ws.Cells.Item(119, 5).Value = 5
$ws.Cells.Item(119, 6).Value = 100112010
$ws.Cells.Item(119, 7).Value = "Achicoria"
$ws.Cells.Item(119, 8).Value = "Sin especificar"
$ws.Cells.Item(119, 9).Value = "Primera"
$ws.Cells.Item(119, 10).Value = 70
$ws.Cells.Item(119, 11).Value = 7000
$ws.Cells.Item(119, 12).Value = 7000
$ws.Cells.Item(119, 13).Value = 7000
$ws.Cells.Item(119, 14).Value = "$/caja 16 unidades"
$ws.Cells.Item(119, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(119, 16).Value = 438
$ws.Cells.Item(119, 17).Value = 16
$ws.Cells.Item(119, 18).Value = "Hortaliza"

# Match the date style (numFmt) used by the other rows in column D
$ws.Cells.Item(119, 4).NumberFormat = $ws.Cells.Item(120, 4).NumberFormat
